$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for the new "Save" column (H), rows 2 through 44
$saveValues = @(1,0,0,1,0,0,1,0,1,0,0,0,0,1,1,0,0,1,1,0,0,0,0,0,0,0,0,0,0,1,0,0,1,1,0,1,0,1,0,0,0,0,1)

# Header cell H1 - copy the existing header formatting (style) from G1, then set the text
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H44 with the save values
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
